$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.944.28'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.331.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.41'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.182'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.581'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.14'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000273'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitcoinCash'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '698.21'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.870.97'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.38'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.963.66'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.119'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.338.19'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.48'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.12'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.893'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Toncoin'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.43'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.91'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '101.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.91'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.69'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.43'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.98'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.54'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.97'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.17%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '574.13'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.01'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.105'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.735.30'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Dai'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.65'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.29'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.60'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +9.91%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.134'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.15'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.60'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0677'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.334'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.31'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0406'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.63'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.129'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '130.26'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.97%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.66'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.36%  '
